$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "15 only address"
$ws.Range("B8").Value = "15 only city"
$ws.Range("C8").Value = " "
$ws.Range("D8").Value = "15 only first"
$ws.Range("E8").Value = "15 only last"
$ws.Range("F8").Value = " "
$ws.Range("G8").Value = "{{address}}"
$ws.Range("H8").Value = "{{address}}"

$ws.Range("A9").Value = "15 address"
$ws.Range("B9").Value = "15 HCM city"
$ws.Range("C9").Value = " "
$ws.Range("D9").Value = "15 silicon first"
$ws.Range("E9").Value = "15 silicon last"
$ws.Range("F9").Value = "{{ip_address}}"
$ws.Range("G9").Value = "'1515"
$ws.Range("G9").Style = "Normal"
$ws.Range("H9").Value = "'151515"
$ws.Range("H9").Style = "Normal"

$ws.Range("A10").Value = "f03 address"
$ws.Range("B10").Value = "03 city"
$ws.Range("C10").Value = " "
$ws.Range("D10").Value = "f03 first"
$ws.Range("E10").Value = "f03 last"
$ws.Range("F10").Value = " "
$ws.Range("G10").Value = "'0303"
$ws.Range("G10").Style = "Normal"
$ws.Range("H10").Value = "'030303"
$ws.Range("H10").Style = "Normal"
